# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right before the "总计" sheet,
#    carrying the single-fund holding snapshot for 2022-Q1 (same layout as
#    the other quarterly sheets).
# 2. Insert a corresponding summary row at the top of the "总计" sheet's
#    data table (row 2), pushing the older quarters down one row, and
#    renumber the leading index column (A).

$wb = $excel.ActiveWorkbook

# A cell that already carries the shared "header / index column" style
# (bold font, box border, centered) used throughout the workbook, so the
# new cells created below match the existing look without inventing a new
# style.
$styleSrc = $wb.Worksheets.Item(5).Range("B1")

# ------------------------------------------------------------------
# 1) Add the new "2022-Q1" sheet just before "总计"
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (row 1) - same bold/bordered style as the other quarter sheets
$styleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data row (row 2) - code/name/size/position columns are stored as text
# (matching the other quarter sheets), the rank column is numeric.
$styleSrc.Copy()
$q1.Range("A2").PasteSpecial(-4122)
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "512040"
$q1.Range("C2").Value = "富国中证价值ETF"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "3.44"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "99.55"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "1.26"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0433"
$q1.Range("H2").Value = 10

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: add a 2022-Q1 row at the top of the table
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$styleSrc.Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.04

# Renumber the index column (A) for the remaining (pre-existing) rows,
# which have each shifted down by one row.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
